# Insert a new weekly price record as row 39 on the active sheet,
# pushing the existing rows 39-49 down to 40-50 (dimension grows to A1:R50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 39 (shifts rows 39:49 -> 40:50).
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new record.
$ws.Cells.Item(39, 1).Value  = 4
$ws.Cells.Item(39, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(39, 3).Value  = "Los Lagos"
$ws.Cells.Item(39, 4).Value  = 44642
$ws.Cells.Item(39, 5).Value  = 10
$ws.Cells.Item(39, 6).Value  = 100112043
$ws.Cells.Item(39, 7).Value  = "Pepino dulce"
$ws.Cells.Item(39, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(39, 9).Value  = "Primera"
$ws.Cells.Item(39, 10).Value = 160
$ws.Cells.Item(39, 11).Value = 17500
$ws.Cells.Item(39, 12).Value = 18000
$ws.Cells.Item(39, 13).Value = 17750
$ws.Cells.Item(39, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 986
$ws.Cells.Item(39, 17).Value = 18
$ws.Cells.Item(39, 18).Value = "Hortaliza"
